$wb = $excel.ActiveWorkbook

# Row 106 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5885067.5
$ws.Range("I106").Value = 6864384.5
$ws.Range("K106").Value = 6864384.5
$ws.Range("M106").Value = -6863753.5

# Row 132 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2044716
$ws.Range("J132").Value = 9094388
$ws.Range("L132").Value = 27283164
$ws.Range("N132").Value = -27288224

# Row 137 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6927.0557
$ws.Range("I137").Value = 8957.388999999999
$ws.Range("J137").Value = 2866.389
$ws.Range("K137").Value = 26872.167
$ws.Range("L137").Value = 8599.167000000001
$ws.Range("M137").Value = -24322.167
$ws.Range("N137").Value = -13699.167

# Row 138 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3743.9146
$ws.Range("J138").Value = 5092.959
$ws.Range("L138").Value = 15278.877
$ws.Range("N138").Value = -25558.877

# Row 22 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1356.5
$ws.Range("I22").Value = 713
$ws.Range("K22").Value = 713
$ws.Range("M22").Value = -414

# Row 28 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 64892068
$ws.Range("I28").Value = 19508.143
$ws.Range("K28").Value = 19508.143
$ws.Range("M28").Value = -19316.143

# Row 32 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4830.63
$ws.Range("I32").Value = 4684.8735
$ws.Range("K32").Value = 4684.8735
$ws.Range("M32").Value = -4397.8735

# Row 61 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5261.4204
$ws.Range("I61").Value = 5407.1963
$ws.Range("K61").Value = 5407.1963
$ws.Range("M61").Value = -5195.1963

# Row 74 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1575.9767
$ws.Range("I74").Value = 921
$ws.Range("K74").Value = 921
$ws.Range("M74").Value = -47

# Row 77 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1575.9767
$ws.Range("I77").Value = 921
$ws.Range("K77").Value = 4605
$ws.Range("M77").Value = -237

# Row 97 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 11118132
$ws.Range("I97").Value = 10042.454
$ws.Range("J97").Value = 28573700
$ws.Range("K97").Value = 10042.454
$ws.Range("L97").Value = 28573700
$ws.Range("M97").Value = -9546.454
$ws.Range("N97").Value = -28574692

# Row 99 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 64892068
$ws.Range("I99").Value = 19508.143
$ws.Range("K99").Value = 19508.143
$ws.Range("M99").Value = -16513.143

# Row 128 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 132 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1602.4565
$ws.Range("I132").Value = 951.3143
$ws.Range("J132").Value = 3674.2727
$ws.Range("K132").Value = 2853.9429
$ws.Range("L132").Value = 11022.8181
$ws.Range("M132").Value = -323.9429
$ws.Range("N132").Value = -16082.8181

# Row 136 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5261.4204
$ws.Range("I136").Value = 5407.1963
$ws.Range("K136").Value = 16221.5889
$ws.Range("M136").Value = -13671.5889

# Row 20 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3469.4243
$ws.Range("I20").Value = 2587.7144
$ws.Range("J20").Value = 5012.4165
$ws.Range("K20").Value = 2587.7144
$ws.Range("L20").Value = 5012.4165
$ws.Range("M20").Value = -2340.7144
$ws.Range("N20").Value = -5506.4165

# Row 94 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2441.3684
$ws.Range("I94").Value = 1389.9
$ws.Range("K94").Value = 1389.9
$ws.Range("M94").Value = -938.9000000000001

# Row 128 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 5000
$ws.Range("I128").Value = 5000
$ws.Range("K128").Value = 15000
$ws.Range("M128").Value = -12510

# Row 7 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 262.78946
$ws.Range("I7").Value = 206.2
$ws.Range("K7").Value = 206.2
$ws.Range("M7").Value = -93.19999999999999

# Row 31 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6028.5815
$ws.Range("I31").Value = 6109.1387
$ws.Range("J31").Value = 5614.2856
$ws.Range("K31").Value = 6109.1387
$ws.Range("L31").Value = 5614.2856
$ws.Range("M31").Value = -5814.1387
$ws.Range("N31").Value = -6204.2856

# Row 34 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6028.5815
$ws.Range("I34").Value = 6109.1387
$ws.Range("J34").Value = 5614.2856
$ws.Range("K34").Value = 6109.1387
$ws.Range("L34").Value = 5614.2856
$ws.Range("M34").Value = -5907.1387
$ws.Range("N34").Value = -6018.2856

# Row 132 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1665.4736
$ws.Range("I132").Value = 1209.6666
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 3628.9998
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -1098.9998
$ws.Range("N132").Value = -15184.25

# Row 141 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 122379.43
$ws.Range("J141").Value = 125983.85
$ws.Range("L141").Value = 125983.85
$ws.Range("N141").Value = -136343.85

# Row 5 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 324549.4
$ws.Range("I5").Value = 998.25
$ws.Range("K5").Value = 2994.75
$ws.Range("M5").Value = -2882.75

# Row 113 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 959
$ws.Range("I113").Value = 713.375
$ws.Range("K113").Value = 2140.125
$ws.Range("M113").Value = 29.875

# Row 122 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4015.4119
$ws.Range("J122").Value = 4972.579
$ws.Range("L122").Value = 44753.211
$ws.Range("N122").Value = -49653.211

# Row 134 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 13046.16
$ws.Range("I134").Value = 28894.375
$ws.Range("J134").Value = 5588.1763
$ws.Range("K134").Value = 86683.125
$ws.Range("L134").Value = 16764.5289
$ws.Range("M134").Value = -81613.125
$ws.Range("N134").Value = -26904.5289

# Row 135 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 324549.4
$ws.Range("I135").Value = 998.25
$ws.Range("K135").Value = 8984.25
$ws.Range("M135").Value = -6449.25

# Row 46 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35333.332
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 35333.332
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 35333.332
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -35645.332

# Row 70 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10248.333
$ws.Range("I70").Value = 8600.429
$ws.Range("J70").Value = 12555.4
$ws.Range("K70").Value = 8600.429
$ws.Range("L70").Value = 12555.4
$ws.Range("M70").Value = -8330.429
$ws.Range("N70").Value = -13095.4

# Row 73 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10248.333
$ws.Range("I73").Value = 8600.429
$ws.Range("J73").Value = 12555.4
$ws.Range("K73").Value = 8600.429
$ws.Range("L73").Value = 12555.4
$ws.Range("M73").Value = -7664.429
$ws.Range("N73").Value = -14427.4

# Row 97 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 6014.8
$ws.Range("I97").Value = 5965.6665
$ws.Range("K97").Value = 5965.6665
$ws.Range("M97").Value = -5469.6665

# Row 102 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6508.2915
$ws.Range("I102").Value = 7529.1577
$ws.Range("J102").Value = 2629
$ws.Range("K102").Value = 7529.1577
$ws.Range("L102").Value = 2629
$ws.Range("M102").Value = -5907.1577
$ws.Range("N102").Value = -5873

# Row 122 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7811.0557
$ws.Range("I122").Value = 5235.7036
$ws.Range("K122").Value = 15707.1108
$ws.Range("M122").Value = -13257.1108

# Row 132 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2121.0588
$ws.Range("I132").Value = 2119.8572
$ws.Range("K132").Value = 6359.571599999999
$ws.Range("M132").Value = -3829.571599999999

# Row 7 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18264.121
$ws.Range("I7").Value = 23946.037
$ws.Range("K7").Value = 23946.037
$ws.Range("M7").Value = -23834.037

# Row 16 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1981.0605
$ws.Range("I16").Value = 1479.1666
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 1479.1666
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -1309.1666
$ws.Range("N16").Value = -7340

# Row 40 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19061.756
$ws.Range("I40").Value = 20895
$ws.Range("J40").Value = 16197.3125
$ws.Range("K40").Value = 20895
$ws.Range("L40").Value = 16197.3125
$ws.Range("M40").Value = -20759
$ws.Range("N40").Value = -16469.3125

# Row 55 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 847.5333000000001
$ws.Range("I55").Value = 388.83334
$ws.Range("K55").Value = 388.83334
$ws.Range("M55").Value = -215.83334

# Row 68 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7950.1113
$ws.Range("I68").Value = 3249
$ws.Range("J68").Value = 9293.286
$ws.Range("K68").Value = 3249
$ws.Range("L68").Value = 9293.286
$ws.Range("M68").Value = -2500
$ws.Range("N68").Value = -10791.286

# Row 71 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7950.1113
$ws.Range("I71").Value = 3249
$ws.Range("J71").Value = 9293.286
$ws.Range("K71").Value = 16245
$ws.Range("L71").Value = 46466.43
$ws.Range("M71").Value = -12501
$ws.Range("N71").Value = -53954.43

# Row 93 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5206.4116
$ws.Range("I93").Value = 5711.5835
$ws.Range("J93").Value = 3994
$ws.Range("K93").Value = 5711.5835
$ws.Range("L93").Value = 3994
$ws.Range("M93").Value = -4463.5835
$ws.Range("N93").Value = -6490

# Row 126 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 18264.121
$ws.Range("I126").Value = 23946.037
$ws.Range("K126").Value = 71838.111
$ws.Range("M126").Value = -69368.111

# Row 62 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 229567.06
$ws.Range("I62").Value = 382973.8
$ws.Range("J62").Value = 10414.571
$ws.Range("K62").Value = 382973.8
$ws.Range("L62").Value = 10414.571
$ws.Range("M62").Value = -382349.8
$ws.Range("N62").Value = -11662.571

# Row 65 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 229567.06
$ws.Range("I65").Value = 382973.8
$ws.Range("J65").Value = 10414.571
$ws.Range("K65").Value = 1914869
$ws.Range("L65").Value = 52072.855
$ws.Range("M65").Value = -1911749
$ws.Range("N65").Value = -58312.855

# Row 122 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4912.4287
$ws.Range("J122").Value = 8202.15
$ws.Range("L122").Value = 24606.45
$ws.Range("N122").Value = -29506.45

# Row 132 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7930.098
$ws.Range("I132").Value = 8824.023999999999
$ws.Range("K132").Value = 26472.072
$ws.Range("M132").Value = -23942.072
